$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 1.668911666666667
$ws.Cells.Item(2, 8).Value = 5.006735
$ws.Cells.Item(2, 9).Value = 0.02533989316516187
$ws.Cells.Item(2, 10).Value = 0.02533989316516188
$ws.Cells.Item(2, 13).Value = 4.114675333333333
$ws.Cells.Item(2, 14).Value = 12.344026
$ws.Cells.Item(2, 15).Value = 0.04794018332925969
$ws.Cells.Item(2, 16).Value = 0.04794018332925969
$ws.Cells.Item(2, 17).Value = 6.867029668345555
$ws.Cells.Item(2, 18).Value = 61.80326701510999
$ws.Cells.Item(2, 19).Value = 0.001214799123881715
$ws.Cells.Item(2, 20).Value = 0.001214799123881715
$ws.Cells.Item(3, 7).Value = 1.668911666666667
$ws.Cells.Item(3, 8).Value = 5.006735
$ws.Cells.Item(3, 9).Value = 0.02533989316516187
$ws.Cells.Item(3, 10).Value = 0.02533989316516188
$ws.Cells.Item(3, 13).Value = 59.62659933333333
$ws.Cells.Item(3, 15).Value = 0.6947109727426806
$ws.Cells.Item(3, 16).Value = 0.6947109727426806
$ws.Cells.Item(3, 17).Value = 99.51152727105888
$ws.Cells.Item(3, 18).Value = 895.6037454395299
$ws.Cells.Item(3, 19).Value = 0.01760390182996521
$ws.Cells.Item(3, 20).Value = 0.01760390182996521
$ws.Cells.Item(4, 7).Value = 1.668911666666667
$ws.Cells.Item(4, 8).Value = 5.006735
$ws.Cells.Item(4, 9).Value = 0.02533989316516187
$ws.Cells.Item(4, 10).Value = 0.02533989316516188
$ws.Cells.Item(4, 13).Value = 21.65107466666666
$ws.Cells.Item(4, 14).Value = 64.95322399999999
$ws.Cells.Item(4, 15).Value = 0.252257202503176
$ws.Cells.Item(4, 16).Value = 0.252257202503176
$ws.Cells.Item(4, 17).Value = 36.1337311070711
$ws.Cells.Item(4, 18).Value = 325.20357996364
$ws.Cells.Item(4, 19).Value = 0.006392170561573084
$ws.Cells.Item(4, 20).Value = 0.006392170561573085
$ws.Cells.Item(5, 7).Value = 1.668911666666667
$ws.Cells.Item(5, 8).Value = 5.006735
$ws.Cells.Item(5, 9).Value = 0.02533989316516187
$ws.Cells.Item(5, 10).Value = 0.02533989316516188
$ws.Cells.Item(5, 13).Value = 0.4370123333333333
$ws.Cells.Item(5, 14).Value = 1.311037
$ws.Cells.Item(5, 15).Value = 0.005091641424883797
$ws.Cells.Item(5, 16).Value = 0.005091641424883797
$ws.Cells.Item(5, 17).Value = 0.7293349815772222
$ws.Cells.Item(5, 18).Value = 6.564014834195
$ws.Cells.Item(5, 19).Value = 0.000129021649741868
$ws.Cells.Item(5, 20).Value = 0.000129021649741868
$ws.Cells.Item(6, 9).Value = 0.409133080066751
$ws.Cells.Item(6, 10).Value = 0.4091330800667511
$ws.Cells.Item(6, 13).Value = 4.114675333333333
$ws.Cells.Item(6, 14).Value = 12.344026
$ws.Cells.Item(6, 15).Value = 0.04794018332925969
$ws.Cells.Item(6, 16).Value = 0.04794018332925969
$ws.Cells.Item(6, 17).Value = 110.8737507616098
$ws.Cells.Item(6, 18).Value = 997.863756854488
$ws.Cells.Item(6, 19).Value = 0.01961391486446473
$ws.Cells.Item(6, 20).Value = 0.01961391486446473
$ws.Cells.Item(7, 9).Value = 0.409133080066751
$ws.Cells.Item(7, 10).Value = 0.4091330800667511
$ws.Cells.Item(7, 13).Value = 59.62659933333333
$ws.Cells.Item(7, 15).Value = 0.6947109727426806
$ws.Cells.Item(7, 16).Value = 0.6947109727426806
$ws.Cells.Item(7, 18).Value = 14460.24718820682
$ws.Cells.Item(7, 19).Value = 0.2842292400343816
$ws.Cells.Item(7, 20).Value = 0.2842292400343817
$ws.Cells.Item(8, 9).Value = 0.409133080066751
$ws.Cells.Item(8, 10).Value = 0.4091330800667511
$ws.Cells.Item(8, 13).Value = 21.65107466666666
$ws.Cells.Item(8, 14).Value = 64.95322399999999
$ws.Cells.Item(8, 15).Value = 0.252257202503176
$ws.Cells.Item(8, 16).Value = 0.252257202503176
$ws.Cells.Item(8, 17).Value = 583.4083279587235
$ws.Cells.Item(8, 18).Value = 5250.674951628512
$ws.Cells.Item(8, 19).Value = 0.1032067662291465
$ws.Cells.Item(8, 20).Value = 0.1032067662291466
$ws.Cells.Item(9, 9).Value = 0.409133080066751
$ws.Cells.Item(9, 10).Value = 0.4091330800667511
$ws.Cells.Item(9, 13).Value = 0.4370123333333333
$ws.Cells.Item(9, 14).Value = 1.311037
$ws.Cells.Item(9, 15).Value = 0.005091641424883797
$ws.Cells.Item(9, 16).Value = 0.005091641424883797
$ws.Cells.Item(9, 17).Value = 11.77570345179511
$ws.Cells.Item(9, 18).Value = 105.981331066156
$ws.Cells.Item(9, 19).Value = 0.002083158938758169
$ws.Cells.Item(9, 20).Value = 0.002083158938758169
$ws.Cells.Item(10, 7).Value = 36.85439666666667
$ws.Cells.Item(10, 8).Value = 110.56319
$ws.Cells.Item(10, 9).Value = 0.5595781327750508
$ws.Cells.Item(10, 10).Value = 0.5595781327750509
$ws.Cells.Item(10, 13).Value = 4.114675333333333
$ws.Cells.Item(10, 14).Value = 12.344026
$ws.Cells.Item(10, 15).Value = 0.04794018332925969
$ws.Cells.Item(10, 16).Value = 0.04794018332925969
$ws.Cells.Item(10, 17).Value = 151.6438768892156
$ws.Cells.Item(10, 18).Value = 1364.79489200294
$ws.Cells.Item(10, 19).Value = 0.02682627827228075
$ws.Cells.Item(10, 20).Value = 0.02682627827228076
$ws.Cells.Item(11, 7).Value = 36.85439666666667
$ws.Cells.Item(11, 8).Value = 110.56319
$ws.Cells.Item(11, 9).Value = 0.5595781327750508
$ws.Cells.Item(11, 10).Value = 0.5595781327750509
$ws.Cells.Item(11, 13).Value = 59.62659933333333
$ws.Cells.Item(11, 15).Value = 0.6947109727426806
$ws.Cells.Item(11, 16).Value = 0.6947109727426806
$ws.Cells.Item(11, 17).Value = 2197.502343715069
$ws.Cells.Item(11, 18).Value = 19777.52109343562
$ws.Cells.Item(11, 19).Value = 0.3887450689456884
$ws.Cells.Item(11, 20).Value = 0.3887450689456884
$ws.Cells.Item(12, 7).Value = 36.85439666666667
$ws.Cells.Item(12, 8).Value = 110.56319
$ws.Cells.Item(12, 9).Value = 0.5595781327750508
$ws.Cells.Item(12, 10).Value = 0.5595781327750509
$ws.Cells.Item(12, 13).Value = 21.65107466666666
$ws.Cells.Item(12, 14).Value = 64.95322399999999
$ws.Cells.Item(12, 15).Value = 0.252257202503176
$ws.Cells.Item(12, 16).Value = 0.252257202503176
$ws.Cells.Item(12, 17).Value = 797.9372940249509
$ws.Cells.Item(12, 18).Value = 7181.43564622456
$ws.Cells.Item(12, 19).Value = 0.1411576143557851
$ws.Cells.Item(12, 20).Value = 0.1411576143557851
$ws.Cells.Item(13, 7).Value = 36.85439666666667
$ws.Cells.Item(13, 8).Value = 110.56319
$ws.Cells.Item(13, 9).Value = 0.5595781327750508
$ws.Cells.Item(13, 10).Value = 0.5595781327750509
$ws.Cells.Item(13, 13).Value = 0.4370123333333333
$ws.Cells.Item(13, 14).Value = 1.311037
$ws.Cells.Item(13, 15).Value = 0.005091641424883797
$ws.Cells.Item(13, 16).Value = 0.005091641424883797
$ws.Cells.Item(13, 17).Value = 16.10582588089222
$ws.Cells.Item(13, 18).Value = 144.95243292803
$ws.Cells.Item(13, 19).Value = 0.002849171201296574
$ws.Cells.Item(13, 20).Value = 0.002849171201296574
$ws.Cells.Item(14, 7).Value = 0.3918003333333333
$ws.Cells.Item(14, 8).Value = 1.175401
$ws.Cells.Item(14, 9).Value = 0.005948893993036266
$ws.Cells.Item(14, 10).Value = 0.005948893993036267
$ws.Cells.Item(14, 13).Value = 4.114675333333333
$ws.Cells.Item(14, 14).Value = 12.344026
$ws.Cells.Item(14, 15).Value = 0.04794018332925969
$ws.Cells.Item(14, 16).Value = 0.04794018332925969
$ws.Cells.Item(14, 17).Value = 1.612131167158444
$ws.Cells.Item(14, 18).Value = 14.509180504426
$ws.Cells.Item(14, 19).Value = 0.0002851910686324903
$ws.Cells.Item(14, 20).Value = 0.0002851910686324903
$ws.Cells.Item(15, 7).Value = 0.3918003333333333
$ws.Cells.Item(15, 8).Value = 1.175401
$ws.Cells.Item(15, 9).Value = 0.005948893993036266
$ws.Cells.Item(15, 10).Value = 0.005948893993036267
$ws.Cells.Item(15, 13).Value = 59.62659933333333
$ws.Cells.Item(15, 15).Value = 0.6947109727426806
$ws.Cells.Item(15, 16).Value = 0.6947109727426806
$ws.Cells.Item(15, 17).Value = 23.36172149433311
$ws.Cells.Item(15, 18).Value = 210.255493448998
$ws.Cells.Item(15, 19).Value = 0.004132761932645314
$ws.Cells.Item(15, 20).Value = 0.004132761932645315
$ws.Cells.Item(16, 7).Value = 0.3918003333333333
$ws.Cells.Item(16, 8).Value = 1.175401
$ws.Cells.Item(16, 9).Value = 0.005948893993036266
$ws.Cells.Item(16, 10).Value = 0.005948893993036267
$ws.Cells.Item(16, 13).Value = 21.65107466666666
$ws.Cells.Item(16, 14).Value = 64.95322399999999
$ws.Cells.Item(16, 15).Value = 0.252257202503176
$ws.Cells.Item(16, 16).Value = 0.252257202503176
$ws.Cells.Item(16, 17).Value = 8.482898271424887
$ws.Cells.Item(16, 18).Value = 76.34608444282398
$ws.Cells.Item(16, 19).Value = 0.001500651356671277
$ws.Cells.Item(16, 20).Value = 0.001500651356671277
$ws.Cells.Item(17, 7).Value = 0.3918003333333333
$ws.Cells.Item(17, 8).Value = 1.175401
$ws.Cells.Item(17, 9).Value = 0.005948893993036266
$ws.Cells.Item(17, 10).Value = 0.005948893993036267
$ws.Cells.Item(17, 13).Value = 0.4370123333333333
$ws.Cells.Item(17, 14).Value = 1.311037
$ws.Cells.Item(17, 15).Value = 0.005091641424883797
$ws.Cells.Item(17, 16).Value = 0.005091641424883797
$ws.Cells.Item(17, 17).Value = 0.1712215778707778
$ws.Cells.Item(17, 18).Value = 1.540994200837
$ws.Cells.Item(17, 19).Value = 0.00003028963508718583
$ws.Cells.Item(17, 20).Value = 0.00003028963508718584

Write-Output "Applied 179 cell updates"
